# Generate Report for Handback
#
# This script applies the "handback" refresh to the localization-status
# report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     for every language row (Overview + each language sheet).
#   - The "Latest Handback DateTime" timestamp is refreshed for zh-cn and de-de.
#   - The stale "handback file is not the latest" error is cleared now that
#     the handback is in sync (Error Detail column).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn (E2) and de-de (F2) status cells -----------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet -------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus                        # Status
$zhcn.Range("K2").Value = "2016-08-19 16:53:26"              # Latest Handback DateTime
$zhcn.Range("P2").Value = ""                                 # Error Detail (cleared)

# --- de-de sheet ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus                         # Status
$dede.Range("K2").Value = "2016-08-19 16:53:33"               # Latest Handback DateTime
$dede.Range("P2").Value = ""                                  # Error Detail (cleared)

# --- Column width refresh (status/error columns grew/shrank with new text) ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668    # E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668    # F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668        # C (Status)
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334       # P (Error Detail)

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668        # C (Status)
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334       # P (Error Detail)
